$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "MAIO 23" row (row 6) figures
$ws.Range("C6").Value = 1100
$ws.Range("E6").Value = 184
$ws.Range("G6").Value = 656
$ws.Range("H6").Value = 1702
$ws.Range("I6").Value = 4084

# Move the active selection to I7 (single cell), matching the author's final cursor position
$ws.Range("I7").Select()
